# Updated symbol list on Mon Dec 19 16:53:43 UTC 2022 with GitHub Actions
#
# Applies the latest cryptocurrency price/label refresh to Sheet1.
# Column D values are stored as text (not numbers) in the workbook, so we
# force each target cell to the "Text" number format before writing the
# new value. This preserves the original "numbers-stored-as-text" layout
# used throughout column D instead of letting Excel auto-convert the
# numeric-looking strings into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Price (column D) updates -------------------------------------------
Set-TextValue "D2"  "245.97"
Set-TextValue "D4"  "5.382"
Set-TextValue "D5"  "0.05660"
Set-TextValue "D6"  "3.381"
Set-TextValue "D7"  "0.8066"
Set-TextValue "D8"  "1.010"
Set-TextValue "D9"  "0.1450"
Set-TextValue "D10" "0.07486"
Set-TextValue "D11" "0.03165"
Set-TextValue "D12" "0.03063"
Set-TextValue "D13" "0.09250"
Set-TextValue "D14" "3.582"
Set-TextValue "D15" "0.001628"
Set-TextValue "D16" "0.04732"
Set-TextValue "D17" "0.0005828"
Set-TextValue "D18" "0.006354"
Set-TextValue "D19" "0.005003"
Set-TextValue "D20" "0.001043"
Set-TextValue "D21" "0.0001502"
Set-TextValue "D22" "0.0003104"
Set-TextValue "D23" "3.768"
Set-TextValue "D24" "6.398"
Set-TextValue "D25" "2.098"
Set-TextValue "D26" "0.3285"
Set-TextValue "D40" "0.04001"
Set-TextValue "D41" "0.006980"
Set-TextValue "D45" "0.00005940"
Set-TextValue "D46" "0.00000000751"
Set-TextValue "D47" "0.0005508"
Set-TextValue "D48" "0.6834"
Set-TextValue "D49" "0.02540"
Set-TextValue "D50" "0.00002103"
Set-TextValue "D51" "0.01011"

# --- Label (column E) tweaks (24h best/worst annotation changes) --------
$ws.Range("E17").Value = "16OneONEWorstin24h"
$ws.Range("E47").Value = "46ACDXExchangeACXT"

# --- Row 42/43 swap: BKEXToken and CEJI traded places in the ranking ----
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003505"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1037"
$ws.Range("E43").Value = "42BKEXTokenBKK"
